$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = [double]"1.221974545399842e-150"
$ws.Cells.Item(2, 4).Value = -0.5516898103439049
$ws.Cells.Item(2, 5).Value = -31.60949773307858
$ws.Cells.Item(3, 1).Value = 0.0268805176360391
$ws.Cells.Item(3, 2).Value = 0.1176918824304137
$ws.Cells.Item(3, 3).Value = 2.427568574858869
$ws.Cells.Item(3, 4).Value = -0.5516898103439049
$ws.Cells.Item(3, 5).Value = -31.60949773307858
$ws.Cells.Item(4, 1).Value = 0.05376103527207819
$ws.Cells.Item(4, 2).Value = 0.2353837648608275
$ws.Cells.Item(4, 3).Value = 2.355137149717738
$ws.Cells.Item(4, 4).Value = -0.5502346828656669
$ws.Cells.Item(4, 5).Value = -31.52612506992203
$ws.Cells.Item(5, 1).Value = 0.08064155290811728
$ws.Cells.Item(5, 2).Value = 0.3531809196109203
$ws.Cells.Item(5, 3).Value = 2.282877057891251
$ws.Cells.Item(5, 4).Value = -0.5473220838147018
$ws.Cells.Item(5, 5).Value = -31.35924543688792
$ws.Cells.Item(6, 1).Value = 0.1075220705441564
$ws.Cells.Item(6, 2).Value = 0.4711880390883038
$ws.Cells.Item(6, 3).Value = 2.210960367960192
$ws.Cells.Item(6, 4).Value = -0.5429473098730622
$ws.Cells.Item(6, 5).Value = -31.10858935370815
$ws.Cells.Item(7, 1).Value = 0.1344025881801955
$ws.Cells.Item(7, 2).Value = 0.5895086475767003
$ws.Cells.Item(7, 3).Value = 2.139560619046592
$ws.Cells.Item(7, 4).Value = -0.5371032680334358
$ws.Cells.Item(7, 5).Value = -30.77375042099969
$ws.Cells.Item(8, 1).Value = 0.1612831058162346
$ws.Cells.Item(8, 2).Value = 0.7082444963242028
$ws.Cells.Item(8, 3).Value = 2.068853556031644
$ws.Cells.Item(8, 4).Value = -0.5297804293754504
$ws.Cells.Item(8, 5).Value = -30.35418267184189
$ws.Cells.Item(9, 1).Value = 0.1881636234522737
$ws.Cells.Item(9, 2).Value = 0.8274949333292863
$ws.Cells.Item(9, 3).Value = 1.999017864497256
$ws.Cells.Item(9, 4).Value = -0.5209667664404288
$ws.Cells.Item(9, 5).Value = -29.84919698361427
$ws.Cells.Item(10, 1).Value = 0.2150441410883128
$ws.Cells.Item(10, 2).Value = 0.9473562389155665
$ws.Cells.Item(10, 3).Value = 1.930235904935595
$ws.Cells.Item(10, 4).Value = -0.5106476733224368
$ws.Cells.Item(10, 5).Value = -29.25795649955083
$ws.Cells.Item(11, 1).Value = 0.2419246587243519
$ws.Cells.Item(11, 2).Value = 1.067920917771034
$ws.Cells.Item(11, 3).Value = 1.862694445442425
$ws.Cells.Item(11, 4).Value = -0.498805867333556
$ws.Cells.Item(11, 5).Value = -28.57947099457521
$ws.Cells.Item(12, 1).Value = 0.268805176360391
$ws.Cells.Item(12, 2).Value = 1.189276937598055
$ws.Cells.Item(12, 3).Value = 1.796585391677201
$ws.Cells.Item(12, 4).Value = -0.4854212708130108
$ws.Cells.Item(12, 5).Value = -27.81259010346249
$ws.Cells.Item(13, 1).Value = 0.2956856939964301
$ws.Cells.Item(13, 2).Value = 1.311506903869867
$ws.Cells.Item(13, 3).Value = 1.73210651231984
$ws.Cells.Item(13, 4).Value = -0.4704708713251782
$ws.Cells.Item(13, 5).Value = -26.95599531077513
$ws.Cells.Item(14, 1).Value = 0.3225662116324691
$ws.Cells.Item(14, 2).Value = 1.434687159409146
$ws.Cells.Item(14, 3).Value = 1.669462157560072
$ws.Cells.Item(14, 4).Value = -0.4539285581212075
$ws.Cells.Item(14, 5).Value = -26.00819058080408
$ws.Cells.Item(15, 1).Value = 0.3494467292685082
$ws.Cells.Item(15, 2).Value = 1.558886796584455
$ws.Cells.Item(15, 3).Value = 1.608863967294098
$ws.Cells.Item(15, 4).Value = -0.4357649323117637
$ws.Cells.Item(15, 5).Value = -24.96749148126805
$ws.Cells.Item(16, 1).Value = 0.3763272469045473
$ws.Cells.Item(16, 2).Value = 1.684166568849727
$ws.Cells.Item(16, 3).Value = 1.550531564642132
$ws.Cells.Item(16, 4).Value = -0.4159470877006611
$ws.Cells.Item(16, 5).Value = -23.8320126260058
$ws.Cells.Item(17, 1).Value = 0.4032077645405864
$ws.Cells.Item(17, 2).Value = 1.810577687117319
$ws.Cells.Item(17, 3).Value = 1.494693229098112
$ws.Cells.Item(17, 4).Value = -0.3944383586441888
$ws.Cells.Item(17, 5).Value = -22.59965322837953
$ws.Cells.Item(18, 1).Value = 0.4300882821766255
$ws.Cells.Item(18, 2).Value = 1.938160485042717
$ws.Cells.Item(18, 3).Value = 1.441586542028435
$ws.Cells.Item(18, 4).Value = -0.3711980306081397
$ws.Cells.Item(18, 5).Value = -21.26808051741436
$ws.Cells.Item(19, 1).Value = 0.4569687998126646
$ws.Cells.Item(19, 2).Value = 2.066942935693423
$ws.Cells.Item(19, 3).Value = 1.391458995286271
$ws.Cells.Item(19, 4).Value = -0.346181008268255
$ws.Cells.Item(19, 5).Value = -19.83471072135447
$ws.Cells.Item(20, 1).Value = 0.4838493174487037
$ws.Cells.Item(20, 2).Value = 2.196939000260382
$ws.Cells.Item(20, 3).Value = 1.344568551322787
$ws.Cells.Item(20, 4).Value = -0.3193374350077572
$ws.Cells.Item(20, 5).Value = -18.29668726647771
$ws.Cells.Item(21, 1).Value = 0.5107298350847428
$ws.Cells.Item(21, 2).Value = 2.32814678743215
$ws.Cells.Item(21, 3).Value = 1.301184140257547
$ws.Cells.Item(21, 4).Value = -0.2906122564671698
$ws.Cells.Item(21, 5).Value = -16.65085577034229
$ws.Cells.Item(22, 1).Value = 0.5376103527207819
$ws.Cells.Item(22, 2).Value = 2.460546499776047
$ws.Cells.Item(22, 3).Value = 1.261586075794204
$ws.Cells.Item(22, 4).Value = -0.2599447193448999
$ws.Cells.Item(22, 5).Value = -14.89373532517545
$ws.Cells.Item(23, 1).Value = 0.564490870356821
$ws.Cells.Item(23, 2).Value = 2.594098140946177
$ws.Cells.Item(23, 3).Value = 1.226066367479603
$ws.Cells.Item(23, 4).Value = -0.2272677948657757
$ws.Cells.Item(23, 5).Value = -13.0214854650539
$ws.Cells.Item(24, 1).Value = 0.5913713879928602
$ws.Cells.Item(24, 2).Value = 2.728738954761222
$ws.Cells.Item(24, 3).Value = 1.194928901410379
$ws.Cells.Item(24, 4).Value = -0.1925075141433495
$ws.Cells.Item(24, 5).Value = -11.02986808496893
$ws.Cells.Item(25, 1).Value = 0.6182519056288992
$ws.Cells.Item(25, 2).Value = 2.864380564171302
$ws.Cells.Item(25, 3).Value = 1.168489454845959
$ws.Cells.Item(25, 4).Value = -0.1555821999492188
$ws.Cells.Item(25, 5).Value = -8.914203424450728
$ws.Cells.Item(26, 1).Value = 0.6451324232649382
$ws.Cells.Item(26, 2).Value = 3.000905774885905
$ws.Cells.Item(26, 3).Value = 1.147075501980531
$ws.Cells.Item(26, 4).Value = -0.1164015760235862
$ws.Cells.Item(26, 5).Value = -6.669319034822683
$ws.Cells.Item(27, 1).Value = 0.6720129409009774
$ws.Cells.Item(27, 2).Value = 3.13816500501282
$ws.Cells.Item(27, 3).Value = 1.131025757965246
$ws.Cells.Item(27, 4).Value = -0.07486573082422021
$ws.Cells.Item(27, 5).Value = -4.289490406390292
$ws.Cells.Item(28, 1).Value = 0.6988934585370165
$ws.Cells.Item(28, 2).Value = 3.27597229855037
$ws.Cells.Item(28, 3).Value = 1.120689395653238
$ws.Cells.Item(28, 4).Value = -0.03086390726137847
$ws.Cells.Item(28, 5).Value = -1.768371625360161
$ws.Cells.Item(29, 1).Value = 0.7257739761730556
$ws.Cells.Item(29, 2).Value = 3.414100877131886
$ws.Cells.Item(29, 3).Value = 1.116424853820902
$ws.Cells.Item(29, 4).Value = 0.01572691683700274
$ws.Cells.Item(29, 5).Value = 0.9010859595134911
$ws.Cells.Item(30, 1).Value = 0.7526544938090947
$ws.Cells.Item(30, 2).Value = 3.552278181283767
$ws.Cells.Item(30, 3).Value = 1.118598135971923
$ws.Cells.Item(30, 4).Value = 0.06504370150367365
$ws.Cells.Item(30, 5).Value = 3.726729580069227
$ws.Cells.Item(31, 1).Value = 0.7795350114451338
$ws.Cells.Item(31, 2).Value = 3.690180350005562
$ws.Cells.Item(31, 3).Value = 1.127580474187301
$ws.Cells.Item(31, 4).Value = 0.1172393258904294
$ws.Cells.Item(31, 5).Value = 6.717318566480448
$ws.Cells.Item(32, 1).Value = 0.8064155290811729
$ws.Cells.Item(32, 2).Value = 3.827426086297172
$ws.Cells.Item(32, 3).Value = 1.143745201458843
$ws.Cells.Item(32, 4).Value = 0.1724849479182087
$ws.Cells.Item(32, 5).Value = 9.882659545247174
$ws.Cells.Item(33, 1).Value = 0.833296046717212
$ws.Cells.Item(33, 2).Value = 3.963569857239248
$ws.Cells.Item(33, 3).Value = 1.167463636728925
$ws.Cells.Item(33, 4).Value = 0.2309728315288771
$ws.Cells.Item(33, 5).Value = 13.23376842879085
$ws.Cells.Item(34, 1).Value = 0.8601765643532511
$ws.Cells.Item(34, 2).Value = 4.098094381738711
$ws.Cells.Item(34, 3).Value = 1.199099737092002
$ws.Cells.Item(34, 4).Value = 0.2929197563761843
$ws.Cells.Item(34, 5).Value = 16.78306577635565
$ws.Cells.Item(35, 1).Value = 0.8870570819892901
$ws.Cells.Item(35, 2).Value = 4.230402369148874
$ws.Cells.Item(35, 3).Value = 1.239003208187407
$ws.Cells.Item(35, 4).Value = 0.3585711591864477
$ws.Cells.Item(35, 5).Value = 20.54461407649705
$ws.Cells.Item(36, 1).Value = 0.9139375996253293
$ws.Cells.Item(36, 2).Value = 4.359807490803462
$ws.Cells.Item(36, 3).Value = 1.287500682632485
$ws.Cells.Item(36, 4).Value = 0.4282062027025618
$ws.Cells.Item(36, 5).Value = 24.53440817618021
$ws.Cells.Item(37, 1).Value = 0.9408181172613683
$ws.Cells.Item(37, 2).Value = 4.485524598854391
$ws.Cells.Item(37, 3).Value = 1.344884471992378
$ws.Cells.Item(37, 4).Value = 0.502144032270546
$ws.Cells.Item(37, 5).Value = 28.7707337567833
$ws.Cells.Item(38, 1).Value = 0.9676986348974075
$ws.Cells.Item(38, 2).Value = 4.606659260010592
$ws.Cells.Item(38, 3).Value = 1.41139826307785
$ws.Cells.Item(38, 4).Value = 0.5807515694460702
$ws.Cells.Item(38, 5).Value = 33.27461387485855
$ws.Cells.Item(39, 1).Value = 0.9945791525334465
$ws.Cells.Item(39, 2).Value = 4.722196757155428
$ws.Cells.Item(39, 3).Value = 1.4872189547886
$ws.Cells.Item(39, 4).Value = 0.6644533181766762
$ws.Cells.Item(39, 5).Value = 38.07037081498677
$ws.Cells.Item(40, 1).Value = 1.021459670169486
$ws.Cells.Item(40, 2).Value = 4.830990846973828
$ws.Cells.Item(40, 3).Value = 1.572433604694444
$ws.Cells.Item(40, 4).Value = 0.7537438401983456
$ws.Cells.Item(40, 5).Value = 43.18634087734836
$ws.Cells.Item(41, 1).Value = 1.048340187805525
$ws.Cells.Item(41, 2).Value = 4.931752774296374
$ws.Cells.Item(41, 3).Value = 1.667010158555513
$ws.Cells.Item(41, 4).Value = 0.8492038206968823
$ws.Cells.Item(41, 5).Value = 48.65579487231566
$ws.Cells.Item(42, 1).Value = 1.075220705441564
$ws.Cells.Item(42, 2).Value = 5.02304137780452
$ws.Cells.Item(42, 3).Value = 1.77076024967606
$ws.Cells.Item(42, 4).Value = 0.9515210387976034
$ws.Cells.Item(42, 5).Value = 54.51813964100653
